$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (pushes existing rows 35..61 down to 36..62)
$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35, 1).Value = 12
$ws.Cells.Item(35, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(35, 3).Value = "Metropolitana"
$ws.Cells.Item(35, 4).Value = 44435
$ws.Cells.Item(35, 5).Value = 13
$ws.Cells.Item(35, 6).Value = 100112002
$ws.Cells.Item(35, 7).Value = "Pimiento"
$ws.Cells.Item(35, 8).Value = "Zafiro verde"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 30
$ws.Cells.Item(35, 11).Value = 35000
$ws.Cells.Item(35, 12).Value = 35000
$ws.Cells.Item(35, 13).Value = 35000
$ws.Cells.Item(35, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(35, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(35, 16).Value = 1944
$ws.Cells.Item(35, 17).Value = 18
$ws.Cells.Item(35, 18).Value = "Hortaliza"

# Append a new row 63 at the end with the latest weekly record
$ws.Cells.Item(63, 1).Value = 12
$ws.Cells.Item(63, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44432
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100112002
$ws.Cells.Item(63, 7).Value = "Pimiento"
$ws.Cells.Item(63, 8).Value = "Zafiro verde"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 20
$ws.Cells.Item(63, 11).Value = 35000
$ws.Cells.Item(63, 12).Value = 35000
$ws.Cells.Item(63, 13).Value = 35000
$ws.Cells.Item(63, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(63, 16).Value = 1944
$ws.Cells.Item(63, 17).Value = 18
$ws.Cells.Item(63, 18).Value = "Hortaliza"

# Give the new date cell the same date number format as the rest of column D
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
